# Updates cryptos list values (price/volume%) and reorders a few coin rows,
# matching the upstream GitHub Actions scrape commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.123.02"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").Value = "2.307.01"
$ws.Range("E3").Value = "  +2.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.74%  "
$ws.Range("E7").Value = "  +2.11%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.507"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.71%  "
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("E13").Value = "  +3.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +16.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").Value = "2.669.05"
$ws.Range("E16").Value = "  +2.92%  "
$ws.Range("D17").Value = "2.291.89"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.809"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.49%  "
$ws.Range("D19").Value = "43.024.44"
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.27%  "
$ws.Range("D21").Value = "0.0₃0906"
$ws.Range("E21").Value = "  +2.60%  "
$ws.Range("E22").Value = "  +3.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.54%  "
$ws.Range("E25").Value = "  +8.83%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.47%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.38%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0699"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("E39").Value = "  +2.43%  "
$ws.Range("E40").Value = "  +6.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.102"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.25%  "
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("D44").Value = "1.995.93"
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("E45").Value = "  +3.78%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.98%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.15%  "
$ws.Range("D50").Value = "2.534.07"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.41%  "
